# Add the ability to import auto pay days and bill mode
$wb = $excel.ActiveWorkbook

# --- Data sheet: new "Auto Pay Days" (J) and "Bill Mode" (K) columns ---
$ws = $wb.Worksheets.Item("Data")

$ws.Range("J1").Value = "Auto Pay Days"
$ws.Range("K1").Value = "Bill Mode"

# Match formatting of the existing header cells (I1)
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:K1").PasteSpecial(-4122) | Out-Null

# Comments describing the new columns
$ws.Range("J1").AddComment("The number of days after invoice generation to run auto pay methods") | Out-Null
$ws.Range("K1").AddComment("Invoice or statement, leave blank to use system default") | Out-Null

# New active cell selection
$ws.Range("K1").Select() | Out-Null

# --- Instructions sheet: resize/reposition the instructions text box ---
$instr = $wb.Worksheets.Item("Instructions")
$shp = $instr.Shapes.Item(1)
$shp.Name = "CustomShape 1"
$shp.Left = 6.066141732283465
$shp.Top = 0.05669291338582677
$shp.Width = 562.4503937007875
$shp.Height = 469.0204724409449

Write-Host "done"
